$d = $word.ActiveDocument

# Surgical plain-text replacement that preserves sibling runs exactly
# (including zero-length "marker" runs like <w:r/>) and never mangles
# straight apostrophes into curly "smart quotes". Word's normal
# Find.Execute(...) replace mode (and setting Range.Text directly)
# both rebuild/coalesce the whole run list of the touched paragraph,
# which drops empty sibling runs and smart-quotes the replacement text
# - so instead we locate the exact character range with Find (no
# replace) and push the raw OOXML for just that range via
# Range.InsertXML, which replaces only that range's contents in place,
# leaving neighboring runs untouched. rPrXml lets us keep the exact
# direct run formatting (bold/italic/...) the original run carried.
function Set-PlainText($range, [string]$newText, [string]$rPrXml = "") {
    $escaped = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:r>' + $rPrXml + '<w:t>' + $escaped + '</w:t></w:r></w:p></w:body>' +
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# Replaces the Nth (1-based) occurrence of $old in the document with
# $new, preserving the run formatting supplied via $rPrXml.
function Replace-Occurrence([string]$old, [string]$new, [int]$occurrence = 1, [string]$rPrXml = "") {
    $text = $d.Content.Text
    $searchFrom = 0
    $idx = -1
    for ($n = 0; $n -lt $occurrence; $n++) {
        $idx = $text.IndexOf($old, $searchFrom)
        if ($idx -lt 0) { break }
        $searchFrom = $idx + 1
    }
    if ($idx -lt 0) {
        throw "Could not find occurrence $occurrence of '$old'"
    }
    $rng = $d.Range($idx, $idx + $old.Length)
    Set-PlainText $rng $new $rPrXml
}

function Replace-AllText([string]$old, [string]$new, [string]$rPrXml = "") {
    while ($true) {
        $text = $d.Content.Text
        $idx = $text.IndexOf($old)
        if ($idx -lt 0) { break }
        $rng = $d.Range($idx, $idx + $old.Length)
        Set-PlainText $rng $new $rPrXml
    }
}

# Title / Heading1 at the very top - plain run, no direct formatting
Replace-Occurrence "Play Beat the Beast: Cerberus' Inferno for Free - Review" "Play Beat the Beast: Cerberus' Inferno for Free" 1

# "What we like" bullet list
Replace-AllText "Wide range of betting options, from 10 cents to 100 euros" "Simple and straightforward gameplay"
Replace-AllText "Exciting bonus game with free spins and 243 ways to win" "Interesting Greek mythology theme"
Replace-AllText "Cerberus symbol acts as both Wild and Scatter for increased winning potential" "Stacked wilds and bonus game with free spins"
Replace-AllText "High RTP of 96.15% and maximum payout of 6,666x stake" "Wide range of betting options"

# "What we don't like" bullet list
Replace-AllText "Limited number of paylines at only 9" "Limited number of paylines"
Replace-AllText "Theme may not appeal to all players" "No progressive jackpot feature"

# Bold "title" run repeated near the end
Replace-Occurrence "Play Beat the Beast: Cerberus' Inferno for Free - Review" "Play Beat the Beast: Cerberus' Inferno for Free" 1 "<w:rPr><w:b/></w:rPr>"

# Meta description (italic run at the end)
Replace-AllText "Read our review of the online slot game Beat the Beast: Cerberus' Inferno, and play for free with exciting features, great payouts, and a Greek mythology theme." "Read our review of Beat the Beast: Cerberus' Inferno and play this slot game for free." "<w:rPr><w:i/></w:rPr>"
